# Solution for "103. Binary Tree Zigzag Level Order Traversal". zen-1
#
# The tracker sheet marks a problem "Finished" by setting column E to "Y"
# (shared string for "N" -> "Y") and hiding that row (the sheet's
# AutoFilter on column E is set to show only "N" rows).
#
# Row 104 is problem 103 "Binary Tree Zigzag Level Order Traversal" - mark
# it solved. Row 83 "Remove Duplicates from Sorted List II" is marked
# solved as well, per the same commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Remove Duplicates from Sorted List II" (row 83) as finished and hide it.
$ws.Range("E83").Value = "Y"
$ws.Rows(83).Hidden = $true

# Mark "103. Binary Tree Zigzag Level Order Traversal" (row 104) as finished and hide it.
$ws.Range("E104").Value = "Y"
$ws.Rows(104).Hidden = $true

# Move the selection/active cell to where the author left off.
$ws.Range("E107").Select()
